$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Handelingen")
$ws.Activate()

# Select entire row 7 (as the user did before deleting it) and delete it,
# shifting all subsequent rows up by one.
$row = $ws.Rows.Item(7)
$row.Select()
$row.Delete()
